{"js": "// Update three Stefan Lyocsa 2025 publication entries (items 21-23) with\n// their correct journal names and DOIs, replacing the placeholder \"SSRN.\"\n// citations.\n\nconst updates = [\n  {\n    find: '21. \"A Fuzzy Framework for Realized Volatility Prediction\" (2025). SSRN.',\n    replace: '21. \"A Fuzzy Framework for Realized Volatility Prediction\" (2025). Journal of Forecasting. DOI: 10.1002/for.70082'\n  },\n  {\n    find: '22. \"Alpha-threshold networks in credit risk models\" (2025). SSRN.',\n    replace: '22. \"Alpha-threshold networks in credit risk models\" (2025). Quantitative Finance. DOI: 10.1080/14697688.2025.2465697'\n  },\n  {\n    find: '23. \"Do hurricanes cause storm on the stock market?\" (2025). SSRN.',\n    replace: '23. \"Do hurricanes cause storm on the stock market?\" (2025). International Review of Financial Analysis. DOI: 10.1016/j.irfa.2024.103816'\n  }\n];\n\nfor (const { find, replace } of updates) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${find}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update three Stefan Lyocsa 2025 publication entries (items 21-23) with\n# their correct journal names and DOIs, replacing the placeholder \"SSRN.\"\n# citations.\n\n$d = $word.ActiveDocument\n\n$updates = @(\n  @{\n    Find    = '21. \"A Fuzzy Framework for Realized Volatility Prediction\" (2025). SSRN.'\n    Replace = '21. \"A Fuzzy Framework for Realized Volatility Prediction\" (2025). Journal of Forecasting. DOI: 10.1002/for.70082'\n  },\n  @{\n    Find    = '22. \"Alpha-threshold networks in credit risk models\" (2025). SSRN.'\n    Replace = '22. \"Alpha-threshold networks in credit risk models\" (2025). Quantitative Finance. DOI: 10.1080/14697688.2025.2465697'\n  },\n  @{\n    Find    = '23. \"Do hurricanes cause storm on the stock market?\" (2025). SSRN.'\n    Replace = '23. \"Do hurricanes cause storm on the stock market?\" (2025). International Review of Financial Analysis. DOI: 10.1016/j.irfa.2024.103816'\n  }\n)\n\nforeach ($update in $updates) {\n  # Use a fresh Range over the whole document body for each search so that\n  # earlier replacements (which change the text length) cannot affect the\n  # search position of later ones.\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Text = $update.Find\n  $find.MatchCase = $true\n  $find.MatchWildcards = $false\n\n  $found = $find.Execute()\n  if (-not $found) {\n    throw \"Could not find text: $($update.Find)\"\n  }\n\n  # Assigning .Text directly on the matched Range replaces its contents\n  # in place without going through AutoCorrect/AutoFormat \"smart quotes\",\n  # so the straight double quotes in the citation are preserved verbatim.\n  $rng.Text = $update.Replace\n}\n\n$d.Save()\n"}
